$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for handback" — the localization round-trip has come
# back from the translators, so the per-language status sheets (zh-cn /
# de-de) get their "target/handback" columns (E..H) filled in, and the
# overview's "Ready for handoff" status flips to "Handed back: in sync
# with en-us" everywhere that text is used (shared string, so setting it
# once per distinct row is enough - Excel re-uses the shared string).
# ---------------------------------------------------------------------------

$zhUrlMd  = "https://github.com/OpenLocalizationTest/oltest/blob/0f6f05446aebfe5746dfe7731ccb86b755d08b91/e2e/512d7bce-4876-42d5-87fd-814ddd8c6112.md"
$zhUrlXlf1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08cd5104cc6b1649ea12eb029443c83216507685/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/512d7bce-4876-42d5-87fd-814ddd8c6112.ca8a3e8c4f571326357d1187970d31981e07bd4a.zh-cn.xlf"
$zhUrlMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/0f6f05446aebfe5746dfe7731ccb86b755d08b91/e2e/79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.md"
$zhUrlXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08cd5104cc6b1649ea12eb029443c83216507685/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.30cbc82b8d020c412954491e6fc4bc651da5c9da.zh-cn.xlf"

$deUrlMd  = "https://github.com/OpenLocalizationTest/oltest/blob/0f6f05446aebfe5746dfe7731ccb86b755d08b91/e2e/512d7bce-4876-42d5-87fd-814ddd8c6112.md"
$deUrlXlf1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/06c5bcb626ca0bd18501ae51c2b388749cda0746/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/512d7bce-4876-42d5-87fd-814ddd8c6112.ca8a3e8c4f571326357d1187970d31981e07bd4a.de-de.xlf"
$deUrlMd2 = "https://github.com/OpenLocalizationTest/oltest/blob/0f6f05446aebfe5746dfe7731ccb86b755d08b91/e2e/79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.md"
$deUrlXlf2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/06c5bcb626ca0bd18501ae51c2b388749cda0746/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.30cbc82b8d020c412954491e6fc4bc651da5c9da.de-de.xlf"

$handedBack = "Handed back: in sync with en-us"

# ---------------------------------------------------------------------------
# Overview sheet - mirrors the per-language Status column
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $handedBack
$ov.Range("C2").Value = $handedBack
$ov.Range("B3").Value = $handedBack
$ov.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Status column now reports the handback instead of "Ready for handoff"
$zh.Range("B2").Value = $handedBack
$zh.Range("B3").Value = $handedBack

# Latest Target File / Latest Handback File (columns E/F) + hyperlinks
$zh.Hyperlinks.Add($zh.Range("E2"), $zhUrlMd, "", "", "512d7bce-4876-42d5-87fd-814ddd8c6112.md")
$zh.Hyperlinks.Add($zh.Range("F2"), $zhUrlXlf1, "", "", "512d7bce-4876-42d5-87fd-814ddd8c6112.ca8a3e8c4f571326357d1187970d31981e07bd4a.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), $zhUrlMd2, "", "", "79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.md")
$zh.Hyperlinks.Add($zh.Range("F3"), $zhUrlXlf2, "", "", "79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.30cbc82b8d020c412954491e6fc4bc651da5c9da.zh-cn.xlf")

# Latest Handback DateTime (G) + Handoff Reason (H)
$zh.Range("G2").Value = "2016-01-09 03:54:07"
$zh.Range("H2").Value = "Include"
$zh.Range("G3").Value = "2016-01-09 03:54:07"
$zh.Range("H3").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $handedBack
$de.Range("B3").Value = $handedBack

$de.Hyperlinks.Add($de.Range("E2"), $deUrlMd, "", "", "512d7bce-4876-42d5-87fd-814ddd8c6112.md")
$de.Hyperlinks.Add($de.Range("F2"), $deUrlXlf1, "", "", "512d7bce-4876-42d5-87fd-814ddd8c6112.ca8a3e8c4f571326357d1187970d31981e07bd4a.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), $deUrlMd2, "", "", "79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.md")
$de.Hyperlinks.Add($de.Range("F3"), $deUrlXlf2, "", "", "79510ec1-f1cd-4c7d-bcb6-48bd2630b61f.30cbc82b8d020c412954491e6fc4bc651da5c9da.de-de.xlf")

$de.Range("G2").Value = "2016-01-09 03:54:24"
$de.Range("H2").Value = "Include"
$de.Range("G3").Value = "2016-01-09 03:54:24"
$de.Range("H3").Value = "Include"

Write-Output "Handback report generated"
